$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly tracker refresh 17/08/2021: new counterfactual (Tracker) values,
# Outturn values carried over per-country, table re-sorted by Tracker ascending.
$data = @(
    @("ESP", -10.84256769219105, -11.19667546947339),
    @("PRT", -9.100051706084356, -10.64625044041106),
    @("ISL", -9.07591809844277, -13.22631292791014),
    @("GBR", -8.564560755108985, -9.852646439789037),
    @("AUT", -7.990040401155996, -9.638008674807697),
    @("CZE", -7.485873550599043, -7.930948186647191),
    @("MEX", -7.430819328173033, -6.978715886918375),
    @("LVA", -7.417625679349705, -7.332354529191843),
    @("GRC", -7.089960034675213, -6.99730348543447),
    @("IDN", -7.029216690540907, -7.889913754437472),
    @("IND", -6.920758258385185, -6.708090798141065),
    @("SVN", -6.691321462594891, -6.866100856477953),
    @("COL", -6.684014578593056, -4.927497657010315),
    @("SVK", -6.602817070781553, -7.108108639465415),
    @("ZAF", -6.543400623016948, -5.57904722784116),
    @("NLD", -6.092857805397223, -5.817335417153624),
    @("DEU", -6.000736438855069, -5.857047898506174),
    @("FRA", -5.837396406909379, -6.195980384081645),
    @("ITA", -5.556722149844551, -7.186499498315612),
    @("POL", -5.382376865440275, -6.537818823294117),
    @("HUN", -5.233089114857991, -6.094651400764151),
    @("JPN", -5.151485355257057, -5.039093521993943),
    @("BEL", -4.748181462766976, -4.718496431588514),
    @("CAN", -4.694540393661317, -3.667707077225879),
    @("CHL", -4.654475237398136, -4.371196753544937),
    @("ROU", -4.344019342569439, -3.537815465093697),
    @("FIN", -4.301980120068427, -3.911233965088456),
    @("NOR", -4.237019737776626, -4.090973536006592),
    @("BGR", -4.024291911737954, -5.021170765135707),
    @("LUX", -4.024050416803138, -0.4768904004482599),
    @("ISR", -3.131566794867058, -5.054740417590542),
    @("DNK", -3.066675649974193, -3.072341303877146),
    @("RUS", -2.987212323246624, -4.763351017441286),
    @("USA", -2.978922762257352, -2.724829140464946),
    @("CHE", -2.866519671286338, -3.204513725473657),
    @("KOR", -2.809799064301033, -1.520234890225636),
    @("LTU", -2.233135578402801, -1.065470275273439),
    @("ARG", -2.20202653678967, -0.5149462598107779),
    @("NZL", -1.914077508116385, -2.394044212023194),
    @("SWE", -1.655265079906953, -2.362646798756485),
    @("AUS", -1.637548339092287, -1.926032824459722),
    @("BRA", -1.161355094468497, -2.555628819729538),
    @("TUR", 1.506455358114889, 4.746288661325582),
    @("EST", 3.06213806776201, 0.80954414990968),
    @("IRL", 4.404226135719469, 10.1997158035124)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

